$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 1
$ws.Range("F2").Value = 2023

# Update existing row 3
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0
$ws.Range("F3").Value = 2975

# Update existing row 4
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0

# Add new row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1

# Copy style from A4 (style index 1: bold, centered, bordered) to A5
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122) # xlPasteFormats
